$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.744.69"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "1.759.78"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "326.57"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.4439"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Value = "0.3736"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "45.46"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "0.07805"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").Value = "1.130"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "21.83"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("D14").Value = "6.202"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "7.388"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "1.759.90"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "91.46"
$ws.Range("E17").Value = "  +13.02%  "
$ws.Range("D18").Value = "0.00001084"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("E19").Value = "  -7.95%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "17.46"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "6.206"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "0.5330"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").Value = "27.770.36"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "11.68"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "2.323"
$ws.Range("E26").Value = "  -3.99%  "
$ws.Range("D27").Value = "20.91"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "153.58"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "2.373"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "1.957.84"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("D31").Value = "129.31"
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("D32").Value = "1.221"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "5.797"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "0.09287"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("E35").Value = "  -9.20%  "
$ws.Range("D36").Value = "12.77"
$ws.Range("E36").Value = "  +5.34%  "
$ws.Range("D37").Value = "0.2201"
$ws.Range("E37").Value = "  -6.46%  "
$ws.Range("D38").Value = "0.02349"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "0.6536"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.113"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.06156"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").Value = "1.201"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "8.052"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").Value = "1.421"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "0.6055"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "3.764"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "126.22"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").Value = "2.003"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "1.152"
$ws.Range("E51").Value = "  -0.63%  "
